# Semi-Auto gun visual application
# Update timesheet: extend totals to row 100 and log new time for LO4 on 2024-02-18 (row 33)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the summed ranges for LO1/utilities (C), LO2 (D) and LO4/world interactions (E)
# totals from row 32 down to row 100 so future entries are picked up automatically.
$ws.Range("C3").Formula = "=SUM(C5:C100)"
$ws.Range("D3").Formula = "=SUM(D5:D100)"
$ws.Range("E3").Formula = "=SUM(E5:E100)"

# Log the time spent on LO4 (world interactions) for the entry in row 33.
$ws.Range("E33").Formula = "=(1/60)*(9+13+21+21)"

# Update the active selection/view to reflect where work left off.
$ws.Range("E34").Select()
